# Generate Report for Handback
# Update timestamps recorded for file 7b2d61c3-682e-4462-bb61-f80f93f3d795.md
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 7b2d61c3... row
$wsOverview.Range("G4").Value = "2016-09-04 02:48:57"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn.Range("H4").Value = "2016-09-04 02:48:52"
$wsZhCn.Range("K4").Value = "2016-09-04 02:49:36"

# de-de sheet: "Correspond Handback DateTime"
$wsDeDe.Range("K4").Value = "2016-09-04 02:49:43"
